$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 13 (pushes existing rows 13-199 down to 14-200) ---
$ws.Rows(13).Insert()

# Carry over the "data row" formatting (with borders) from row 12 onto the
# three newly available rows (13, 14, 15) so they match the surrounding
# table look instead of the stripped blank-row formatting Insert() leaves.
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 13 : new timesheet entry ---
$ws.Range("A13").Value2 = 45964
$ws.Range("B13").Value2 = "Research"
$ws.Range("C13").Value2 = 0.55972222222222223
$ws.Range("D13").Value2 = 0.65833333333333333
$ws.Range("E13").Formula = "=D13-C13"
$ws.Range("F13").Value2 = "Researched more about the Dart language syntax and using Flutter"

# --- Fix typo in the F9 description cell ---
$ws.Range("F9").Value2 = "Check-in: Group meeting with Tyler, discussed what was needed for the pre-alpha build, planning more research for app development"

# --- Row 14 : new timesheet entry ---
$ws.Range("A14").Value2 = 45965
$ws.Range("B14").Value2 = "Check-in"
$ws.Range("C14").Value2 = 0.47222222222222221
$ws.Range("D14").Value2 = 0.49305555555555558
$ws.Range("E14").Formula = "=D14-C14"
$ws.Range("F14").Value2 = "Check-in: Group meeting with Tyler, discussed communication between frontend and backend,  tasks to complete, and next steps for group coding session"

# --- Row 15 : new timesheet entry ---
$ws.Range("A15").Value2 = 45967
$ws.Range("B15").Value2 = "Frontend Dev"
$ws.Range("C15").Value2 = 0.4201388888888889
$ws.Range("D15").Value2 = 0.61527777777777781
$ws.Range("E15").Formula = "=D15-C15"
$ws.Range("F15").Value2 = "Implemented addPlantView screen, reformatted home screen(route clickable button), and established temporary connection between front/backend. Can get data from the backend server"

# --- Restore the view state: scrolled back to A1, F15 selected ---
$ws.Range("A1").Select()
$ws.Range("F15").Select()
